$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-9 (header row 1 stays as-is).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.
$data = @(
    @("M2","Cd28","Cd80","ECs", 3,1,7.543204,22.629612,0.9718881576768906,0.9718881576768906,2,0.6666666666666666,0.539733,1.619199,0.0338608224813917,0.0338608224813917,4.071316124532,36.64184512078801,0.03290893237886402,0.03290893237886402),
    @("M2","Cd28","Cd80","FAPs", 3,1,7.543204,22.629612,0.9718881576768906,0.9718881576768906,3,1,3.713472666666666,11.140418,0.2329693362375475,0.2329693362375475,28.01148187309066,252.103336857816,0.2264201389911181,0.2264201389911181),
    @("M2","Cd28","Cd80","M2",   3,1,7.543204,22.629612,0.9718881576768906,0.9718881576768906,3,1,10.63672166666667,31.910165,0.6673079914308979,0.6673079914308979,80.23496142288667,722.1146528059801,0.6485487343948417,0.6485487343948417),
    @("M2","Cd28","Cd80","sCs",  3,1,7.543204,22.629612,0.9718881576768906,0.9718881576768906,3,1,1.049821333333333,3.149464,0.06586184985016284,0.06586184985016284,7.919016480885333,71.27114832796801,0.06401035191206676,0.06401035191206676),
    @("sCs","Cd28","Cd80","ECs", 1,0.3333333333333333,0.218187,0.6545609999999999,0.02811184232310935,0.02811184232310934,2,0.6666666666666666,0.539733,1.619199,0.0338608224813917,0.0338608224813917,0.117762724071,1.059864516639,0.0009518901025276795,0.0009518901025276794),
    @("sCs","Cd28","Cd80","FAPs",1,0.3333333333333333,0.218187,0.6545609999999999,0.02811184232310935,0.02811184232310934,3,1,3.713472666666666,11.140418,0.2329693362375475,0.2329693362375475,0.8102314607219999,7.292083146497998,0.00654919724642938,0.006549197246429379),
    @("sCs","Cd28","Cd80","M2",  1,0.3333333333333333,0.218187,0.6545609999999999,0.02811184232310935,0.02811184232310934,3,1,10.63672166666667,31.910165,0.6673079914308979,0.6673079914308979,2.320794390285,20.887149512565,0.0187592570360562,0.0187592570360562),
    @("sCs","Cd28","Cd80","sCs", 1,0.3333333333333333,0.218187,0.6545609999999999,0.02811184232310935,0.02811184232310934,3,1,1.049821333333333,3.149464,0.06586184985016284,0.06586184985016284,0.229057367256,2.061516305304,0.001851497938096081,0.001851497938096081)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
